# Loan RBI, Variable Instalments
#
# 1. Insert a new (blank) column before column N ("Late") on the
#    "Repayment schedule" sheet, shifting the existing N/O/P columns one
#    place to the right (-> O/P/Q). The newly inserted column inherits the
#    width of the column immediately to its left (column M = 11 chars).
# 2. Make "Repayment schedule" the active sheet/tab and select cell K13
#    there (this is also what puts tabSelected="1" on that sheet and
#    removes it from whichever sheet was previously active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns.Item(14).Insert()

# ColumnWidth is expressed in "characters" and Excel stores a slightly
# different number internally (characters + ~0.8333 for the default
# Calibri 11 gridline/padding offset on this workbook). 10.1666... set
# here round-trips to the stored width="11" that column M already uses,
# matching the formatting Excel applies when inserting a column (it
# copies the format of the column to the left).
$ws.Columns.Item(14).ColumnWidth = 10.166666666666666

$ws.Activate()
$ws.Range("K13").Select()
